{"js": "// Target change (see diff): inside word/numbering.xml, the <w:nsid> GUID\n// stamped on four <w:abstractNum> definitions (abstractNumId 990, 991,\n// 99721, 99722) is replaced by a newly generated GUID:\n//   990   : 478df75a -> 48ea52dd\n//   991   : 46f0d268 -> 75b35c4d\n//   99721 : acdde4bc -> a833eaef\n//   99722 : 25a05475 -> e922b72e\n// Nothing else in the package changes: multiLevelType, every <w:lvl>\n// (numFmt/lvlText/indents/...), the <w:num> -> <w:abstractNum> mapping, and\n// every paragraph's numPr/ilvl stay byte-identical. <w:nsid> is Word's\n// internal \"list GUID\" bookkeeping value (akin to the rsid pool in\n// settings.xml) - it is written by Word whenever it (re)persists a list\n// definition, but it is never surfaced as a scriptable property: neither\n// the Word JS API (Word.List / Word.ListLevel / Word.Paragraph) nor the\n// Word object model exposes \"nsid\" or any OOXML-part editor capable of\n// reaching it, so it cannot be toggled from a task pane / script lab add-in\n// either. InsertOoxml/getOoxml only ever touch body (\"document.xml\")\n// content - any other <pkg:part> included in a WordOpenXML payload is\n// ignored - and there is no API that lets an add-in rewrite word/numbering.xml\n// directly, so the GUID swap itself is out of reach of the Office.js object\n// model (exactly as in real Word).\n//\n// What *is* addressable is the document content that actually consumes\n// those list definitions, so walk it and confirm it is untouched: resolve\n// the two paragraphs that are list items and read back the list id Office.js\n// does expose (Word.List.id, which is the <w:num numId>, not the <w:nsid>\n// GUID) together with their list level, without writing anything back.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  paragraph.load(\"isListItem\");\n}\nawait context.sync();\n\nconst listItems = paragraphs.items.filter((p) => p.isListItem);\nfor (const paragraph of listItems) {\n  paragraph.list.load(\"id\");\n}\nawait context.sync();\n\n// No body/content mutation is made: the only bytes the diff touches\n// (word/numbering.xml's <w:nsid> values) have no corresponding\n// context.document.* API to set, so there is nothing safe to write here\n// without inventing an unrelated, unrequested edit.\n", "ps1": "# Target change (see diff): inside word/numbering.xml, the <w:nsid> GUID\n# stamped on four <w:abstractNum> definitions (abstractNumId 990, 991,\n# 99721, 99722) is replaced by a newly generated GUID:\n#   990   : 478df75a -> 48ea52dd\n#   991   : 46f0d268 -> 75b35c4d\n#   99721 : acdde4bc -> a833eaef\n#   99722 : 25a05475 -> e922b72e\n# Nothing else in the package changes: multiLevelType, every <w:lvl>\n# (numFmt/lvlText/indents/...), the <w:num> -> <w:abstractNum> mapping, and\n# every paragraph's ListFormat stay identical. <w:nsid> is Word's internal\n# \"list GUID\" bookkeeping value (akin to the rsid pool in settings.xml) -\n# Word stamps/rewrites it whenever it persists a list definition, but it is\n# not a documented/scriptable property anywhere in the Word object model:\n# List/ListTemplate/ListLevel expose ListID (the <w:num numId>), not the\n# <w:nsid> GUID, and there is no AbstractNum object at all. Range.WordOpenXML\n# (and Document.WordOpenXML) can *read* the whole package, including\n# word/numbering.xml as a <pkg:part>, but both are read-only - assigning\n# back to them raises \"is a read-only property\" - and Range.InsertXML only\n# ever replaces the body content of the target range, not other package\n# parts. So the GUID swap itself cannot be produced from this object model\n# (exactly as in real Word, where no VBA/COM call reaches <w:nsid> either).\n#\n# What *is* addressable is the document content that actually consumes\n# those list definitions, so walk it and confirm it is untouched: find the\n# paragraphs attached to a list and read back the identifiers COM does\n# expose (ListFormat.List.ListID, which is the <w:num numId>, and\n# ListFormat.ListLevelNumber), without writing anything back.\n$d = $word.ActiveDocument\n\nforeach ($paragraph in $d.Paragraphs) {\n    $range = $paragraph.Range\n    if ($range.ListFormat.ListType -ne 0) {\n        $listId = $range.ListFormat.List.ListID\n        $level = $range.ListFormat.ListLevelNumber\n    }\n}\n\n# No content mutation is made: the only bytes the diff touches\n# (word/numbering.xml's <w:nsid> values) have no corresponding\n# ActiveDocument.* property or method to set, so there is nothing safe to\n# write here without inventing an unrelated, unrequested edit.\n"}
